# Insert a new weekly price record at row 849, shifting the existing
# rows 849-896 down to 850-897 (matches the "Fruta / hortaliza, semanal"
# weekly refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(849).Insert()

$ws.Range("A849").Value = 6
$ws.Range("B849").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C849").Value = "Metropolitana"
$ws.Range("D849").Value = 44706
$ws.Range("E849").Value = 13
$ws.Range("F849").Value = 100112031
$ws.Range("G849").Value = "Poroto verde"
$ws.Range("H849").Value = "Magnum"
$ws.Range("I849").Value = "Primera"
$ws.Range("J849").Value = 250
$ws.Range("K849").Value = 21000
$ws.Range("L849").Value = 23000
$ws.Range("M849").Value = 21800
$ws.Range("N849").Value = "$/saco 25 kilos"
$ws.Range("O849").Value = "Región Metropolitana"
$ws.Range("P849").Value = 872
$ws.Range("Q849").Value = 25
$ws.Range("R849").Value = "Hortaliza"
